$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The template used to ship with 5 blank sanction rows (rows 3-6); trim
# down to just the 2 extra sample rows that are actually filled in below.
$ws.Rows("5:6").Delete()

# Row 2 was a half-filled sample (EmpUnqID/ShiftCode/TPAHours) - wipe the
# columns that are not part of the new sample data before repopulating.
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()

# --- Sample data rows (entered as text via leading apostrophe so Excel
# doesn't coerce the ids/dates/times into numbers) ---

# Row 2
$ws.Range("B2").Value = "'2020-08-26"
$ws.Range("A2").Value = "'104019"
$ws.Range("C2").Value = "'09:00"

# Row 3
$ws.Range("B3").Value = "'2020-08-25"
$ws.Range("A3").Value = "'104019"
$ws.Range("C3").Value = "'09:00"

# Row 4
$ws.Range("A4").Value = "'112244"
$ws.Range("B4").Value = "'2020-08-26"
$ws.Range("C4").Value = "'08:00"
$ws.Range("D4").Value = "'20:00"

# OutTime filled in afterwards for the first two rows
$ws.Range("D2").Value = "'18:00"
$ws.Range("D3").Value = "'18:00"

$ws.Range("F4").Value = 4

# New "Reason" column header, bold like the other headers
$ws.Range("G1").Value = "Reason"
$ws.Range("G1").Font.Bold = $true

# G4 keeps the same (quote-prefixed) style as the other data cells but
# stays empty
$ws.Range("G4").Value = "'x"
$ws.Range("G4").ClearContents()

# Column B needs to grow to fit the new (longer) sample values, and
# column G is widened to comfortably fit future "Reason" text.
$ws.Columns("B").ColumnWidth = 11.3
$ws.Columns("G").ColumnWidth = 14.75

$ws.Range("G2").Select()
